# Fix full/empty circuit diagram labels: "set empty" and "set full" were
# backwards, and nudge/resize their text boxes to match the corrected
# label widths. Also refresh the cached "today" date shown by the
# datetimeFigureOut footer fields on the slide master and every slide
# layout (PowerPoint recalculates these on save).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Swap the "set empty" / "set full" textbox labels on slide 1, and
#    correct their stored position/size so the wording fits.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    if (-not $sh.TextFrame.HasText) { continue }

    $t = $sh.TextFrame.TextRange.Text

    if ($t -eq "set empty") {
        $sh.Left = 123.289685
        $sh.Top = 292.856614
        $sh.Width = 124.379292
        $sh.Height = 29.081260
        $sh.TextFrame.TextRange.Text = "set full"
    }
    elseif ($t -eq "set full") {
        $sh.Left = 743.537953
        $sh.Top = 151.944253
        $sh.Width = 90.352917
        $sh.Height = 29.081260
        $sh.TextFrame.TextRange.Text = "set empty"
    }
}

# ---------------------------------------------------------------------
# 2) Refresh the cached date text on the "Date Placeholder" field shape
#    everywhere it appears: the slide master and all slide layouts.
# ---------------------------------------------------------------------
function Update-DatePlaceholders {
    param($shapes, $newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name.StartsWith("Date Placeholder")) {
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.HasText) {
                    $sh.TextFrame.TextRange.Text = $newText
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "11/23/2020"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholders $lay.Shapes "11/23/2020"
}
